$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header anchor words (unchanged text, indices shift internally)
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"

# Update anchor-word block (A-H), rows 3-6; row 7 A-H cleared entirely
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.8529411764705882
$ws.Range("C3").Value = 29
$ws.Range("D3").Value = 29
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 5

$ws.Range("A4").Value = "crisis"
$ws.Range("B4").Value = 0.589041095890411
$ws.Range("C4").Value = 172
$ws.Range("D4").Value = 172
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 120

$ws.Range("A5").Value = "panic"
$ws.Range("B5").Value = 0.1705426356589147
$ws.Range("C5").Value = 88
$ws.Range("D5").Value = 88
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 428

$ws.Range("A6").Value = "sc"
$ws.Range("B6").Value = 0.1693121693121693
$ws.Range("C6").Value = 32
$ws.Range("D6").Value = 32
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 157

# Row 7 A:H no longer has data (sheet shrank by one anchor-word row)
$ws.Range("A7:H7").Clear() | Out-Null

# Rows 28-29 are brand new; give J28/J29 the same label formatting (bold + border)
# used by the existing J3:J27 label cells, by copying the format from J27.
$ws.Range("J27").Copy() | Out-Null
$ws.Range("J28").PasteSpecial(-4122) | Out-Null
$ws.Range("J29").PasteSpecial(-4122) | Out-Null

# Update keyword confidence block (J-Q), rows 3-29
$ws.Range("J3").Value = "love"
$ws.Range("K3").Value = 0.9782608695652174
$ws.Range("L3").Value = 45
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 1

$ws.Range("J4").Value = "happy"
$ws.Range("K4").Value = 0.9615384615384616
$ws.Range("L4").Value = 25
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 1

$ws.Range("J5").Value = "interesting"
$ws.Range("K5").Value = 0.9393939393939394
$ws.Range("L5").Value = 31
$ws.Range("M5").Value = 31
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 2

$ws.Range("J6").Value = "best"
$ws.Range("K6").Value = 0.9152542372881356
$ws.Range("L6").Value = 54
$ws.Range("M6").Value = 54
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 5

$ws.Range("J7").Value = "great"
$ws.Range("K7").Value = 0.8571428571428571
$ws.Range("L7").Value = 96
$ws.Range("M7").Value = 96
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 16

$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8292682926829268
$ws.Range("L8").Value = 68
$ws.Range("M8").Value = 68
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 14

$ws.Range("J9").Value = "special"
$ws.Range("K9").Value = 0.7777777777777778
$ws.Range("L9").Value = 28
$ws.Range("M9").Value = 28
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 8

$ws.Range("J10").Value = "thank"
$ws.Range("K10").Value = 0.7734375
$ws.Range("L10").Value = 99
$ws.Range("M10").Value = 99
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 29

$ws.Range("J11").Value = "confidence"
$ws.Range("K11").Value = 0.75
$ws.Range("L11").Value = 27
$ws.Range("M11").Value = 27
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 9

$ws.Range("J12").Value = "free"
$ws.Range("K12").Value = 0.7416666666666667
$ws.Range("L12").Value = 89
$ws.Range("M12").Value = 89
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 31

$ws.Range("J13").Value = "positive"
$ws.Range("K13").Value = 0.7413793103448276
$ws.Range("L13").Value = 43
$ws.Range("M13").Value = 43
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 15

$ws.Range("J14").Value = "safe"
$ws.Range("K14").Value = 0.7394366197183099
$ws.Range("L14").Value = 105
$ws.Range("M14").Value = 105
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 37

$ws.Range("J15").Value = "support"
$ws.Range("K15").Value = 0.7075471698113207
$ws.Range("L15").Value = 75
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 31

$ws.Range("J16").Value = "safety"
$ws.Range("K16").Value = 0.6666666666666666
$ws.Range("L16").Value = 34
$ws.Range("M16").Value = 34
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 17

$ws.Range("J17").Value = "good"
$ws.Range("K17").Value = 0.6625
$ws.Range("L17").Value = 106
$ws.Range("M17").Value = 106
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 54

$ws.Range("J18").Value = "better"
$ws.Range("K18").Value = 0.6507936507936508
$ws.Range("L18").Value = 41
$ws.Range("M18").Value = 41
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 22

$ws.Range("J19").Value = "relief"
$ws.Range("K19").Value = 0.62
$ws.Range("L19").Value = 31
$ws.Range("M19").Value = 31
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 19

$ws.Range("J20").Value = "fresh"
$ws.Range("K20").Value = 0.5833333333333334
$ws.Range("L20").Value = 28
$ws.Range("M20").Value = 28
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 20

$ws.Range("J21").Value = "well"
$ws.Range("K21").Value = 0.574468085106383
$ws.Range("L21").Value = 54
$ws.Range("M21").Value = 54
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 40

$ws.Range("J22").Value = "care"
$ws.Range("K22").Value = 0.4943820224719101
$ws.Range("L22").Value = 44
$ws.Range("M22").Value = 44
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 45

$ws.Range("J23").Value = "like"
$ws.Range("K23").Value = 0.4676470588235294
$ws.Range("L23").Value = 159
$ws.Range("M23").Value = 159
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 181

$ws.Range("J24").Value = "hand"
$ws.Range("K24").Value = 0.4621409921671018
$ws.Range("L24").Value = 177
$ws.Range("M24").Value = 177
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 206

$ws.Range("J25").Value = "help"
$ws.Range("K25").Value = 0.4203389830508474
$ws.Range("L25").Value = 124
$ws.Range("M25").Value = 124
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 171

$ws.Range("J26").Value = "protect"
$ws.Range("K26").Value = 0.3698630136986301
$ws.Range("L26").Value = 27
$ws.Range("M26").Value = 27
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 46

$ws.Range("J27").Value = "increase"
$ws.Range("K27").Value = 0.358974358974359
$ws.Range("L27").Value = 28
$ws.Range("M27").Value = 28
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 50

$ws.Range("J28").Value = "please"
$ws.Range("K28").Value = 0.3138075313807531
$ws.Range("L28").Value = 75
$ws.Range("M28").Value = 75
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 164

$ws.Range("J29").Value = "store"
$ws.Range("K29").Value = 0.0302013422818792
$ws.Range("L29").Value = 27
$ws.Range("M29").Value = 27
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 867
